$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new (blank) column before column U. This shifts the existing
#    U/V/W/X columns (rows 7-11 and 18) one column to the right (-> V/W/X/Y),
#    matching the "slideshow 3" block added alongside the existing two blocks.
# ---------------------------------------------------------------------------
$ws.Columns("U").Insert()

# ---------------------------------------------------------------------------
# 2. New header row (row 6) for the shifted block (V6:Y6) and the brand new
#    third block (AB6:AE6) - same header text pattern as M6:P6.
# ---------------------------------------------------------------------------
$ws.Range("V6").Value = "*,"
$ws.Range("W6").Value = "*{}"
$ws.Range("X6").Value = "translate"
$ws.Range("Y6").Value = "img"

$ws.Range("AB6").Value = "*,"
$ws.Range("AC6").Value = "*{}"
$ws.Range("AD6").Value = "translate"
$ws.Range("AE6").Value = "img"

# ---------------------------------------------------------------------------
# 3. New data column V (rows 7-11) for the shifted block, plus the brand new
#    third block AB:AE (rows 7-9).
# ---------------------------------------------------------------------------
$ws.Range("V7").Value = 0
$ws.Range("V8").Value = 20
$ws.Range("V9").Value = 40
$ws.Range("V10").Value = 60
$ws.Range("V11").Value = 80

$ws.Range("AB7").Value = 0
$ws.Range("AC7").Value = 18.666666666666664
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 1

$ws.Range("AB8").Value = 33.3333
$ws.Range("AC8").Formula = "=AC7+AB8"
$ws.Range("AD8").Value = -33.333
$ws.Range("AE8").Value = 2

$ws.Range("AB9").Value = 66.6667
$ws.Range("AC9").Formula = "=AC7+AB9"
$ws.Range("AD9").Value = -66.667
$ws.Range("AE9").Value = 3

# ---------------------------------------------------------------------------
# 4. New labels/formulas around row 18-20.
# ---------------------------------------------------------------------------
$ws.Range("N18").Value = "int 7"
$ws.Range("Z18").Formula = "=11.2/20"

$ws.Range("N19").Value = "step 12.5"

$ws.Range("AC20").Formula = "=100/3"
$ws.Range("AD20").Formula = "=AC20*Z18"

# ---------------------------------------------------------------------------
# 5. View state - scroll / selection (best effort; cosmetic only).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 20
$win.ScrollRow = 1
$ws.Range("AC10").Select()
